# Applies the "tilføjer sælger table på normal former" edit:
#  1. Append ", sælgernavn" to the Låntilbud (1NF) list, right after ", bilpris".
#  2. Append ", sælgerid" to the Låntilbud (3NF) list, right after "BilId".
#  3. Move the "_GoBack" bookmark from the start of the "Bil: ..." paragraph to
#     the end of the "Låntilbud: ... BilId, sælgerid" paragraph (collapsed,
#     right before that paragraph's mark).
#  4. Add a new paragraph "Sælger: id, navn" right after "Bil: BilId, pris, navn".

$d = $word.ActiveDocument

# --- 1. ", sælgernavn" after ", bilpris" -----------------------------------
$r = $d.Content.Duplicate
$null = $r.Find.Execute("bilpris", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)   # wdCollapseEnd
$r.InsertAfter(", sælgernavn")

# --- 2. ", sælgerid" after "BilId" in the 3NF Låntilbud paragraph ----------
$r = $d.Content.Duplicate
$null = $r.Find.Execute("Renteset, BilId", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)   # wdCollapseEnd
$r.InsertAfter(", sælgerid")

# --- 3. Move the _GoBack bookmark ------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the paragraph that now ends in "...BilId, sælgerid" and work out a
# range that sits right before its paragraph mark (its Range.End - 1).
# A bookmark collapsed exactly at "paragraph end - 1" confuses the engine, so
# we anchor it next to a temporary marker character and then remove the
# marker, which keeps the bookmark pinned at the right spot.
$r = $d.Content.Duplicate
$null = $r.Find.Execute("BilId, sælgerid", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)   # wdCollapseEnd, right before the paragraph mark
$r.InsertAfter("~")

$r2 = $d.Content.Duplicate
$null = $r2.Find.Execute("sælgerid~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerStart = $r2.End - 1

$bmRange = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$rMarker = $d.Range($markerStart, $markerStart + 1)
$rMarker.Delete()

# --- 4. New paragraph "Sælger: id, navn" -----------------------------------
$r = $d.Content.Duplicate
$null = $r.Find.Execute("Bil: BilId, pris, navn", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)   # wdCollapseEnd
$r.InsertParagraphAfter()

$r3 = $d.Content.Duplicate
$null = $r3.Find.Execute("Bil: BilId, pris, navn", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Collapse(0)
$r3.MoveStart(1, 1)   # wdCharacter, move past the paragraph mark we just inserted
$r3.InsertAfter("Sælger: id, navn")
